$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = 'Última actualización: 06:53:31'
$ws1.Range("A3").Value = 'Total filas: 63'
$ws1.Cells.Item(6, 1).Value = '04:03:00'
$ws1.Cells.Item(6, 2).Value = '04:03'
$ws1.Cells.Item(6, 3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(6, 4).Value = 0
$ws1.Cells.Item(6, 5).Value = 'LP1912'
$ws1.Cells.Item(7, 1).Value = '04:37:19'
$ws1.Cells.Item(7, 2).Value = '04:46'
$ws1.Cells.Item(7, 3).Value = '215A_EL PATO'
$ws1.Cells.Item(7, 4).Value = 9
$ws1.Cells.Item(7, 5).Value = 'LP1912'
$ws1.Cells.Item(8, 1).Value = '04:52:25'
$ws1.Cells.Item(8, 2).Value = '04:53'
$ws1.Cells.Item(8, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(8, 4).Value = 1
$ws1.Cells.Item(8, 5).Value = 'LP1912'
$ws1.Cells.Item(9, 1).Value = '04:52:25'
$ws1.Cells.Item(9, 2).Value = '05:16'
$ws1.Cells.Item(9, 3).Value = '17_ROMERO'
$ws1.Cells.Item(9, 4).Value = 24
$ws1.Cells.Item(9, 5).Value = 'LP1912'
$ws1.Cells.Item(10, 1).Value = '05:20:30'
$ws1.Cells.Item(10, 2).Value = '05:20'
$ws1.Cells.Item(10, 3).Value = '17_ROMERO'
$ws1.Cells.Item(10, 4).Value = 0
$ws1.Cells.Item(10, 5).Value = 'LP1912'
$ws1.Cells.Item(11, 1).Value = '04:52:25'
$ws1.Cells.Item(11, 2).Value = '05:22'
$ws1.Cells.Item(11, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(11, 4).Value = 30
$ws1.Cells.Item(11, 5).Value = 'LP1912'
$ws1.Cells.Item(12, 1).Value = '05:20:30'
$ws1.Cells.Item(12, 2).Value = '05:26'
$ws1.Cells.Item(12, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(12, 4).Value = 6
$ws1.Cells.Item(12, 5).Value = 'LP1912'
$ws1.Cells.Item(13, 1).Value = '05:20:30'
$ws1.Cells.Item(13, 2).Value = '05:34'
$ws1.Cells.Item(13, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(13, 4).Value = 14
$ws1.Cells.Item(13, 5).Value = 'LP1912'
$ws1.Cells.Item(14, 1).Value = '04:03:00'
$ws1.Cells.Item(14, 2).Value = '05:35'
$ws1.Cells.Item(14, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(14, 4).Value = 92
$ws1.Cells.Item(14, 5).Value = 'LP1912'
$ws1.Cells.Item(15, 1).Value = '04:03:00'
$ws1.Cells.Item(15, 2).Value = '05:41'
$ws1.Cells.Item(15, 3).Value = '14_ABASTO'
$ws1.Cells.Item(15, 4).Value = 98
$ws1.Cells.Item(15, 5).Value = 'LP1912'
$ws1.Cells.Item(16, 1).Value = '05:20:30'
$ws1.Cells.Item(16, 2).Value = '05:46'
$ws1.Cells.Item(16, 3).Value = '15_ABASTO'
$ws1.Cells.Item(16, 4).Value = 26
$ws1.Cells.Item(16, 5).Value = 'LP1912'
$ws1.Cells.Item(17, 1).Value = '05:54:55'
$ws1.Cells.Item(17, 2).Value = '05:54'
$ws1.Cells.Item(17, 3).Value = '10_OLMOS'
$ws1.Cells.Item(17, 4).Value = 0
$ws1.Cells.Item(17, 5).Value = 'LP1912'
$ws1.Cells.Item(18, 1).Value = '05:54:55'
$ws1.Cells.Item(18, 2).Value = '05:55'
$ws1.Cells.Item(18, 3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(18, 4).Value = 1
$ws1.Cells.Item(18, 5).Value = 'LP1912'
$ws1.Cells.Item(19, 1).Value = '05:20:30'
$ws1.Cells.Item(19, 2).Value = '06:04'
$ws1.Cells.Item(19, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(19, 4).Value = 44
$ws1.Cells.Item(19, 5).Value = 'LP1912'
$ws1.Cells.Item(20, 1).Value = '05:54:55'
$ws1.Cells.Item(20, 2).Value = '06:11'
$ws1.Cells.Item(20, 3).Value = '215A_EL PATO'
$ws1.Cells.Item(20, 4).Value = 17
$ws1.Cells.Item(20, 5).Value = 'LP1912'
$ws1.Cells.Item(21, 1).Value = '05:54:55'
$ws1.Cells.Item(21, 2).Value = '06:13'
$ws1.Cells.Item(21, 3).Value = '225_HARAS DEL SUR'
$ws1.Cells.Item(21, 4).Value = 19
$ws1.Cells.Item(21, 5).Value = 'LP1912'
$ws1.Cells.Item(22, 1).Value = '05:20:30'
$ws1.Cells.Item(22, 2).Value = '06:14'
$ws1.Cells.Item(22, 3).Value = '225_HARAS DEL SUR'
$ws1.Cells.Item(22, 4).Value = 54
$ws1.Cells.Item(22, 5).Value = 'LP1912'
$ws1.Cells.Item(23, 1).Value = '05:54:55'
$ws1.Cells.Item(23, 2).Value = '06:20'
$ws1.Cells.Item(23, 3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(23, 4).Value = 26
$ws1.Cells.Item(23, 5).Value = 'LP1912'
$ws1.Cells.Item(24, 1).Value = '05:20:30'
$ws1.Cells.Item(24, 2).Value = '06:21'
$ws1.Cells.Item(24, 3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(24, 4).Value = 61
$ws1.Cells.Item(24, 5).Value = 'LP1912'
$ws1.Cells.Item(25, 1).Value = '05:54:55'
$ws1.Cells.Item(25, 2).Value = '06:26'
$ws1.Cells.Item(25, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(25, 4).Value = 32
$ws1.Cells.Item(25, 5).Value = 'LP1912'
$ws1.Cells.Item(26, 1).Value = '06:24:16'
$ws1.Cells.Item(26, 2).Value = '06:27'
$ws1.Cells.Item(26, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(26, 4).Value = 3
$ws1.Cells.Item(26, 5).Value = 'LP1912'
$ws1.Cells.Item(27, 1).Value = '05:54:55'
$ws1.Cells.Item(27, 2).Value = '06:29'
$ws1.Cells.Item(27, 3).Value = '86_EST CHICA-ESC AGRARIA'
$ws1.Cells.Item(27, 4).Value = 35
$ws1.Cells.Item(27, 5).Value = 'LP1912'
$ws1.Cells.Item(28, 1).Value = '06:24:16'
$ws1.Cells.Item(28, 2).Value = '06:30'
$ws1.Cells.Item(28, 3).Value = '86_EST CHICA-ESC AGRARIA'
$ws1.Cells.Item(28, 4).Value = 6
$ws1.Cells.Item(28, 5).Value = 'LP1912'
$ws1.Cells.Item(29, 1).Value = '06:24:16'
$ws1.Cells.Item(29, 2).Value = '06:31'
$ws1.Cells.Item(29, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(29, 4).Value = 7
$ws1.Cells.Item(29, 5).Value = 'LP1912'
$ws1.Cells.Item(30, 1).Value = '05:54:55'
$ws1.Cells.Item(30, 2).Value = '06:43'
$ws1.Cells.Item(30, 3).Value = '225_C ROCA-H SUR'
$ws1.Cells.Item(30, 4).Value = 49
$ws1.Cells.Item(30, 5).Value = 'LP1912'
$ws1.Cells.Item(31, 1).Value = '06:24:16'
$ws1.Cells.Item(31, 2).Value = '06:44'
$ws1.Cells.Item(31, 3).Value = '225_C ROCA-H SUR'
$ws1.Cells.Item(31, 4).Value = 20
$ws1.Cells.Item(31, 5).Value = 'LP1912'
$ws1.Cells.Item(32, 1).Value = '06:24:16'
$ws1.Cells.Item(32, 2).Value = '06:46'
$ws1.Cells.Item(32, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(32, 4).Value = 22
$ws1.Cells.Item(32, 5).Value = 'LP1912'
$ws1.Cells.Item(33, 1).Value = '06:53:31'
$ws1.Cells.Item(33, 2).Value = '06:56'
$ws1.Cells.Item(33, 3).Value = '14_ABASTO'
$ws1.Cells.Item(33, 4).Value = 3
$ws1.Cells.Item(33, 5).Value = 'LP1912'
$ws1.Cells.Item(34, 1).Value = '05:54:55'
$ws1.Cells.Item(34, 2).Value = '06:59'
$ws1.Cells.Item(34, 3).Value = '14_ABASTO'
$ws1.Cells.Item(34, 4).Value = 65
$ws1.Cells.Item(34, 5).Value = 'LP1912'
$ws1.Cells.Item(35, 1).Value = '06:24:16'
$ws1.Cells.Item(35, 2).Value = '07:00'
$ws1.Cells.Item(35, 3).Value = '14_ABASTO'
$ws1.Cells.Item(35, 4).Value = 36
$ws1.Cells.Item(35, 5).Value = 'LP1912'
$ws1.Cells.Item(36, 1).Value = '06:53:31'
$ws1.Cells.Item(36, 2).Value = '07:01'
$ws1.Cells.Item(36, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(36, 4).Value = 8
$ws1.Cells.Item(36, 5).Value = 'LP1912'
$ws1.Cells.Item(37, 1).Value = '06:53:31'
$ws1.Cells.Item(37, 2).Value = '07:04'
$ws1.Cells.Item(37, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(37, 4).Value = 11
$ws1.Cells.Item(37, 5).Value = 'LP1912'
$ws1.Cells.Item(38, 1).Value = '06:53:31'
$ws1.Cells.Item(38, 2).Value = '07:05'
$ws1.Cells.Item(38, 3).Value = '15_ABASTO'
$ws1.Cells.Item(38, 4).Value = 12
$ws1.Cells.Item(38, 5).Value = 'LP1912'
$ws1.Cells.Item(39, 1).Value = '05:54:55'
$ws1.Cells.Item(39, 2).Value = '07:06'
$ws1.Cells.Item(39, 3).Value = '225_GOMEZ'
$ws1.Cells.Item(39, 4).Value = 72
$ws1.Cells.Item(39, 5).Value = 'LP1912'
$ws1.Cells.Item(40, 1).Value = '06:53:31'
$ws1.Cells.Item(40, 2).Value = '07:07'
$ws1.Cells.Item(40, 3).Value = '225_GOMEZ'
$ws1.Cells.Item(40, 4).Value = 14
$ws1.Cells.Item(40, 5).Value = 'LP1912'
$ws1.Cells.Item(41, 1).Value = '06:53:31'
$ws1.Cells.Item(41, 2).Value = '07:11'
$ws1.Cells.Item(41, 3).Value = '215A_EL PATO'
$ws1.Cells.Item(41, 4).Value = 18
$ws1.Cells.Item(41, 5).Value = 'LP1912'
$ws1.Cells.Item(42, 1).Value = '06:53:31'
$ws1.Cells.Item(42, 2).Value = '07:15'
$ws1.Cells.Item(42, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(42, 4).Value = 22
$ws1.Cells.Item(42, 5).Value = 'LP1912'
$ws1.Cells.Item(43, 1).Value = '06:24:16'
$ws1.Cells.Item(43, 2).Value = '07:16'
$ws1.Cells.Item(43, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(43, 4).Value = 52
$ws1.Cells.Item(43, 5).Value = 'LP1912'
$ws1.Cells.Item(44, 1).Value = '06:53:31'
$ws1.Cells.Item(44, 2).Value = '07:16'
$ws1.Cells.Item(44, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(44, 4).Value = 23
$ws1.Cells.Item(44, 5).Value = 'LP1912'
$ws1.Cells.Item(45, 1).Value = '05:54:55'
$ws1.Cells.Item(45, 2).Value = '07:20'
$ws1.Cells.Item(45, 3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(45, 4).Value = 86
$ws1.Cells.Item(45, 5).Value = 'LP1912'
$ws1.Cells.Item(46, 1).Value = '06:53:31'
$ws1.Cells.Item(46, 2).Value = '07:21'
$ws1.Cells.Item(46, 3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(46, 4).Value = 28
$ws1.Cells.Item(46, 5).Value = 'LP1912'
$ws1.Cells.Item(47, 1).Value = '05:54:55'
$ws1.Cells.Item(47, 2).Value = '07:22'
$ws1.Cells.Item(47, 3).Value = '10_OLMOS'
$ws1.Cells.Item(47, 4).Value = 88
$ws1.Cells.Item(47, 5).Value = 'LP1912'
$ws1.Cells.Item(48, 1).Value = '06:53:31'
$ws1.Cells.Item(48, 2).Value = '07:23'
$ws1.Cells.Item(48, 3).Value = '10_OLMOS'
$ws1.Cells.Item(48, 4).Value = 30
$ws1.Cells.Item(48, 5).Value = 'LP1912'
$ws1.Cells.Item(49, 1).Value = '06:53:31'
$ws1.Cells.Item(49, 2).Value = '07:31'
$ws1.Cells.Item(49, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(49, 4).Value = 38
$ws1.Cells.Item(49, 5).Value = 'LP1912'
$ws1.Cells.Item(50, 1).Value = '05:54:55'
$ws1.Cells.Item(50, 2).Value = '07:31'
$ws1.Cells.Item(50, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(50, 4).Value = 97
$ws1.Cells.Item(50, 5).Value = 'LP1912'
$ws1.Cells.Item(51, 1).Value = '06:53:31'
$ws1.Cells.Item(51, 2).Value = '07:32'
$ws1.Cells.Item(51, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(51, 4).Value = 39
$ws1.Cells.Item(51, 5).Value = 'LP1912'
$ws1.Cells.Item(52, 1).Value = '06:24:16'
$ws1.Cells.Item(52, 2).Value = '07:32'
$ws1.Cells.Item(52, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(52, 4).Value = 68
$ws1.Cells.Item(52, 5).Value = 'LP1912'
$ws1.Cells.Item(53, 1).Value = '06:53:31'
$ws1.Cells.Item(53, 2).Value = '07:36'
$ws1.Cells.Item(53, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(53, 4).Value = 43
$ws1.Cells.Item(53, 5).Value = 'LP1912'
$ws1.Cells.Item(54, 1).Value = '06:24:16'
$ws1.Cells.Item(54, 2).Value = '07:37'
$ws1.Cells.Item(54, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(54, 4).Value = 73
$ws1.Cells.Item(54, 5).Value = 'LP1912'
$ws1.Cells.Item(55, 1).Value = '06:24:16'
$ws1.Cells.Item(55, 2).Value = '07:39'
$ws1.Cells.Item(55, 3).Value = '10_OLMOS'
$ws1.Cells.Item(55, 4).Value = 75
$ws1.Cells.Item(55, 5).Value = 'LP1912'
$ws1.Cells.Item(56, 1).Value = '06:53:31'
$ws1.Cells.Item(56, 2).Value = '07:47'
$ws1.Cells.Item(56, 3).Value = '14_ABASTO'
$ws1.Cells.Item(56, 4).Value = 54
$ws1.Cells.Item(56, 5).Value = 'LP1912'
$ws1.Cells.Item(57, 1).Value = '06:24:16'
$ws1.Cells.Item(57, 2).Value = '07:48'
$ws1.Cells.Item(57, 3).Value = '14_ABASTO'
$ws1.Cells.Item(57, 4).Value = 84
$ws1.Cells.Item(57, 5).Value = 'LP1912'
$ws1.Cells.Item(58, 1).Value = '06:53:31'
$ws1.Cells.Item(58, 2).Value = '07:51'
$ws1.Cells.Item(58, 3).Value = '215D_EL PATO'
$ws1.Cells.Item(58, 4).Value = 58
$ws1.Cells.Item(58, 5).Value = 'LP1912'
$ws1.Cells.Item(59, 1).Value = '06:53:31'
$ws1.Cells.Item(59, 2).Value = '07:55'
$ws1.Cells.Item(59, 3).Value = '10_OLMOS'
$ws1.Cells.Item(59, 4).Value = 62
$ws1.Cells.Item(59, 5).Value = 'LP1912'
$ws1.Cells.Item(60, 1).Value = '06:24:16'
$ws1.Cells.Item(60, 2).Value = '08:00'
$ws1.Cells.Item(60, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(60, 4).Value = 96
$ws1.Cells.Item(60, 5).Value = 'LP1912'
$ws1.Cells.Item(61, 1).Value = '06:53:31'
$ws1.Cells.Item(61, 2).Value = '08:05'
$ws1.Cells.Item(61, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(61, 4).Value = 72
$ws1.Cells.Item(61, 5).Value = 'LP1912'
$ws1.Cells.Item(62, 1).Value = '06:53:31'
$ws1.Cells.Item(62, 2).Value = '08:12'
$ws1.Cells.Item(62, 3).Value = '15_ABASTO'
$ws1.Cells.Item(62, 4).Value = 79
$ws1.Cells.Item(62, 5).Value = 'LP1912'
$ws1.Cells.Item(63, 1).Value = '06:53:31'
$ws1.Cells.Item(63, 2).Value = '08:21'
$ws1.Cells.Item(63, 3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(63, 4).Value = 88
$ws1.Cells.Item(63, 5).Value = 'LP1912'
$ws1.Cells.Item(64, 1).Value = '06:53:31'
$ws1.Cells.Item(64, 2).Value = '08:22'
$ws1.Cells.Item(64, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(64, 4).Value = 89
$ws1.Cells.Item(64, 5).Value = 'LP1912'
$ws1.Cells.Item(65, 1).Value = '06:24:16'
$ws1.Cells.Item(65, 2).Value = '08:23'
$ws1.Cells.Item(65, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(65, 4).Value = 119
$ws1.Cells.Item(65, 5).Value = 'LP1912'
$ws1.Cells.Item(66, 1).Value = '06:53:31'
$ws1.Cells.Item(66, 2).Value = '08:23'
$ws1.Cells.Item(66, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(66, 4).Value = 90
$ws1.Cells.Item(66, 5).Value = 'LP1912'
$ws1.Cells.Item(67, 1).Value = '06:53:31'
$ws1.Cells.Item(67, 2).Value = '08:27'
$ws1.Cells.Item(67, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(67, 4).Value = 94
$ws1.Cells.Item(67, 5).Value = 'LP1912'
$ws1.Cells.Item(68, 1).Value = '06:53:31'
$ws1.Cells.Item(68, 2).Value = '08:42'
$ws1.Cells.Item(68, 3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(68, 4).Value = 109
$ws1.Cells.Item(68, 5).Value = 'LP1912'

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = 'Última actualización: 06:53:31'
$ws2.Cells.Item(6, 1).Value = '04:37:19'
$ws2.Cells.Item(6, 2).Value = '04:46'
$ws2.Cells.Item(6, 3).Value = '215A_EL PATO'
$ws2.Cells.Item(6, 4).Value = 9
$ws2.Cells.Item(6, 5).Value = 'LP1912'
$ws2.Cells.Item(7, 1).Value = '05:20:30'
$ws2.Cells.Item(7, 2).Value = '05:34'
$ws2.Cells.Item(7, 3).Value = '215B_EL PATO'
$ws2.Cells.Item(7, 4).Value = 14
$ws2.Cells.Item(7, 5).Value = 'LP1912'
$ws2.Cells.Item(8, 1).Value = '04:03:00'
$ws2.Cells.Item(8, 2).Value = '05:35'
$ws2.Cells.Item(8, 3).Value = '215B_EL PATO'
$ws2.Cells.Item(8, 4).Value = 92
$ws2.Cells.Item(8, 5).Value = 'LP1912'
$ws2.Cells.Item(9, 1).Value = '05:54:55'
$ws2.Cells.Item(9, 2).Value = '06:11'
$ws2.Cells.Item(9, 3).Value = '215A_EL PATO'
$ws2.Cells.Item(9, 4).Value = 17
$ws2.Cells.Item(9, 5).Value = 'LP1912'
$ws2.Cells.Item(10, 1).Value = '06:24:16'
$ws2.Cells.Item(10, 2).Value = '06:46'
$ws2.Cells.Item(10, 3).Value = '215C_EL PATO'
$ws2.Cells.Item(10, 4).Value = 22
$ws2.Cells.Item(10, 5).Value = 'LP1912'
$ws2.Cells.Item(11, 1).Value = '06:53:31'
$ws2.Cells.Item(11, 2).Value = '07:11'
$ws2.Cells.Item(11, 3).Value = '215A_EL PATO'
$ws2.Cells.Item(11, 4).Value = 18
$ws2.Cells.Item(11, 5).Value = 'LP1912'
$ws2.Cells.Item(12, 1).Value = '06:53:31'
$ws2.Cells.Item(12, 2).Value = '07:51'
$ws2.Cells.Item(12, 3).Value = '215D_EL PATO'
$ws2.Cells.Item(12, 4).Value = 58
$ws2.Cells.Item(12, 5).Value = 'LP1912'
$ws2.Cells.Item(13, 1).Value = '06:53:31'
$ws2.Cells.Item(13, 2).Value = '08:23'
$ws2.Cells.Item(13, 3).Value = '215B_EL PATO'
$ws2.Cells.Item(13, 4).Value = 90
$ws2.Cells.Item(13, 5).Value = 'LP1912'

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = 'Última actualización: 06:53:31'
$ws3.Range("A3").Value = 'Total filas: 11'
$ws3.Cells.Item(6, 1).Value = '05:20:30'
$ws3.Cells.Item(6, 2).Value = '05:44'
$ws3.Cells.Item(6, 3).Value = '215A_LA PLATA'
$ws3.Cells.Item(6, 4).Value = 24
$ws3.Cells.Item(6, 5).Value = 'L6173'
$ws3.Cells.Item(7, 1).Value = '05:54:55'
$ws3.Cells.Item(7, 2).Value = '06:09'
$ws3.Cells.Item(7, 3).Value = '215A_LA PLATA'
$ws3.Cells.Item(7, 4).Value = 15
$ws3.Cells.Item(7, 5).Value = 'L6173'
$ws3.Cells.Item(8, 1).Value = '05:20:30'
$ws3.Cells.Item(8, 2).Value = '06:10'
$ws3.Cells.Item(8, 3).Value = '215A_LA PLATA'
$ws3.Cells.Item(8, 4).Value = 50
$ws3.Cells.Item(8, 5).Value = 'L6173'
$ws3.Cells.Item(9, 1).Value = '05:54:55'
$ws3.Cells.Item(9, 2).Value = '06:32'
$ws3.Cells.Item(9, 3).Value = '215C_LA PLATA'
$ws3.Cells.Item(9, 4).Value = 38
$ws3.Cells.Item(9, 5).Value = 'L6203'
$ws3.Cells.Item(10, 1).Value = '06:24:16'
$ws3.Cells.Item(10, 2).Value = '06:33'
$ws3.Cells.Item(10, 3).Value = '215C_LA PLATA'
$ws3.Cells.Item(10, 4).Value = 9
$ws3.Cells.Item(10, 5).Value = 'L6203'
$ws3.Cells.Item(11, 1).Value = '05:54:55'
$ws3.Cells.Item(11, 2).Value = '06:59'
$ws3.Cells.Item(11, 3).Value = '215B_LP-P MOR-1 Y 57'
$ws3.Cells.Item(11, 4).Value = 65
$ws3.Cells.Item(11, 5).Value = 'L6173'
$ws3.Cells.Item(12, 1).Value = '06:53:31'
$ws3.Cells.Item(12, 2).Value = '07:00'
$ws3.Cells.Item(12, 3).Value = '215B_LP-P MOR-1 Y 57'
$ws3.Cells.Item(12, 4).Value = 7
$ws3.Cells.Item(12, 5).Value = 'L6173'
$ws3.Cells.Item(13, 1).Value = '05:54:55'
$ws3.Cells.Item(13, 2).Value = '07:34'
$ws3.Cells.Item(13, 3).Value = '215A_LA PLATA'
$ws3.Cells.Item(13, 4).Value = 100
$ws3.Cells.Item(13, 5).Value = 'L6173'
$ws3.Cells.Item(14, 1).Value = '06:53:31'
$ws3.Cells.Item(14, 2).Value = '07:35'
$ws3.Cells.Item(14, 3).Value = '215A_LA PLATA'
$ws3.Cells.Item(14, 4).Value = 42
$ws3.Cells.Item(14, 5).Value = 'L6173'
$ws3.Cells.Item(15, 1).Value = '06:53:31'
$ws3.Cells.Item(15, 2).Value = '08:07'
$ws3.Cells.Item(15, 3).Value = '215C_LA PLATA'
$ws3.Cells.Item(15, 4).Value = 74
$ws3.Cells.Item(15, 5).Value = 'L6203'
$ws3.Cells.Item(16, 1).Value = '06:53:31'
$ws3.Cells.Item(16, 2).Value = '08:35'
$ws3.Cells.Item(16, 3).Value = '215A_LA PLATA'
$ws3.Cells.Item(16, 4).Value = 102
$ws3.Cells.Item(16, 5).Value = 'L6173'
